$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.8824145423368052
$ws.Range("E2").Value = 0.8824145423368052

# Row 3
$ws.Range("D3").Value = 0.8769092692153508
$ws.Range("E3").Value = 0.8769092692153508

# Row 4
$ws.Range("D4").Value = 0.6703762919267402
$ws.Range("E4").Value = 0.6703762919267402

# Row 5
$ws.Range("C5").Value = $true
$ws.Range("D5").Value = 0.2752038298194512
$ws.Range("E5").Value = 0.2752038298194512

# Row 6
$ws.Range("C6").Value = $true
$ws.Range("D6").Value = 0.07237308967002483
$ws.Range("E6").Value = 0.07237308967002483

# Row 7
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = 0.03241780870271954
$ws.Range("E7").Value = 0.9675821912972804

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = 0.3240371710350318
$ws.Range("E8").Value = 0.6759628289649682

# Row 9
$ws.Range("D9").Value = 0.5446593746588304
$ws.Range("E9").Value = 0.4553406253411696

# Row 10
$ws.Range("D10").Value = 0.6639029842498251
$ws.Range("E10").Value = 0.3360970157501749

# Row 11
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.3342725366030656
$ws.Range("E11").Value = 0.6657274633969344
$ws.Range("F11").Value = 1.241117358207703
$ws.Range("G11").Value = 0.4

# Row 12
$ws.Range("D12").Value = 0.9717628801031348
$ws.Range("E12").Value = 0.9717628801031348

# Row 13
$ws.Range("D13").Value = 0.9403129537294537
$ws.Range("E13").Value = 0.9403129537294537

# Row 14
$ws.Range("D14").Value = 0.8192659231572341
$ws.Range("E14").Value = 0.8192659231572341

# Row 15
$ws.Range("C15").Value = $true
$ws.Range("D15").Value = 0.2601880179357907
$ws.Range("E15").Value = 0.2601880179357907

# Row 16
$ws.Range("C16").Value = $true
$ws.Range("D16").Value = 0.03002036616483027
$ws.Range("E16").Value = 0.03002036616483027

# Row 17
$ws.Range("C17").Value = $false
$ws.Range("D17").Value = 0.00677333948388859
$ws.Range("E17").Value = 0.9932266605161114

# Row 18
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = 0.2648327109808429
$ws.Range("E18").Value = 0.7351672890191571

# Row 19
$ws.Range("D19").Value = 0.6341751954120193
$ws.Range("E19").Value = 0.3658248045879807

# Row 20
$ws.Range("D20").Value = 0.6949513306309474
$ws.Range("E20").Value = 0.3050486693690526

# Row 21
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = 0.3260708628823657
$ws.Range("E21").Value = 0.6739291371176344
$ws.Range("F21").Value = 1.669172644615173
$ws.Range("G21").Value = 0.4
